$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'40.098.55"
$ws.Range("E2").Value = "'  +1.20%  "
$ws.Range("D3").Value = "'2.218.15"
$ws.Range("E3").Value = "'  +0.13%  "
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("D5").Value = "'290.17"
$ws.Range("E5").Value = "'  -2.84%  "
$ws.Range("D6").Value = "'88.35"
$ws.Range("E6").Value = "'  +5.63%  "
$ws.Range("E7").Value = "'  +0.21%  "
$ws.Range("E8").Value = "'  -0.02%  "
$ws.Range("E9").Value = "'  +1.39%  "
$ws.Range("D10").Value = "'30.79"
$ws.Range("E10").Value = "'  +4.10%  "
$ws.Range("D11").Value = "'0.0782"
$ws.Range("E11").Value = "'  +0.26%  "
$ws.Range("D12").Value = "'47.93"
$ws.Range("E12").Value = "'  +3.75%  "
$ws.Range("E13").Value = "'  +2.70%  "
$ws.Range("D14").Value = "'6.47"
$ws.Range("E14").Value = "'  +3.15%  "
$ws.Range("D15").Value = "'2.559.06"
$ws.Range("E15").Value = "'  +0.05%  "
$ws.Range("D16").Value = "'14.04"
$ws.Range("E16").Value = "'  -0.44%  "
$ws.Range("D17").Value = "'2.214.47"
$ws.Range("E17").Value = "'  -0.24%  "
$ws.Range("D18").Value = "'0.729"
$ws.Range("E18").Value = "'  +1.60%  "
$ws.Range("D19").Value = "'40.028.74"
$ws.Range("E19").Value = "'  +1.20%  "
$ws.Range("D20").Value = "'11.93"
$ws.Range("E20").Value = "'  +14.75%  "
$ws.Range("D21").Value = "'0.0₃0887"
$ws.Range("E21").Value = "'  +1.05%  "
$ws.Range("D22").Value = "'5.82"
$ws.Range("E22").Value = "'  +1.56%  "
$ws.Range("D23").Value = "'65.60"
$ws.Range("E23").Value = "'  +0.99%  "
$ws.Range("D24").Value = "'235.64"
$ws.Range("E24").Value = "'  +1.34%  "
$ws.Range("E25").Value = "'  +0.11%  "
$ws.Range("D26").Value = "'2.46"
$ws.Range("E26").Value = "'  +1.65%  "
$ws.Range("D27").Value = "'1.84"
$ws.Range("E27").Value = "'  +0.61%  "
$ws.Range("D28").Value = "'22.65"
$ws.Range("E28").Value = "'  -0.26%  "
$ws.Range("E29").Value = "'  +1.13%  "
$ws.Range("E30").Value = "'  +1.01%  "
$ws.Range("D31").Value = "'153.26"
$ws.Range("E31").Value = "'  +2.70%  "
$ws.Range("D32").Value = "'32.19"
$ws.Range("E32").Value = "'  -0.25%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "'  -0.05%  "
$ws.Range("D34").Value = "'4.96"
$ws.Range("E34").Value = "'  +2.79%  "
$ws.Range("E35").Value = "'  +2.91%  "
$ws.Range("E36").Value = "'  -0.05%  "
$ws.Range("D37").Value = "'2.84"
$ws.Range("E37").Value = "'  +7.09%  "
$ws.Range("B38").Value = "'Stellar"
$ws.Range("C38").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "'0.112"
$ws.Range("E38").Value = "'  +0.63%  "
$ws.Range("B39").Value = "'Celestia"
$ws.Range("C39").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'15.96"
$ws.Range("E39").Value = "'  -0.99%  "
$ws.Range("D40").Value = "'0.0999"
$ws.Range("E40").Value = "'  +2.83%  "
$ws.Range("E41").Value = "'  +3.38%  "
$ws.Range("B42").Value = "'Maker"
$ws.Range("C42").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'2.096.62"
$ws.Range("E42").Value = "'  +8.87%  "
$ws.Range("B43").Value = "'RenderToken"
$ws.Range("C43").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'3.84"
$ws.Range("E43").Value = "'  +5.04%  "
$ws.Range("D44").Value = "'2.18"
$ws.Range("E44").Value = "'  +1.74%  "
$ws.Range("E45").Value = "'  +1.33%  "
$ws.Range("D46").Value = "'9.88"
$ws.Range("E46").Value = "'  +7.30%  "
$ws.Range("D47").Value = "'17.71"
$ws.Range("E47").Value = "'  +7.91%  "
$ws.Range("D48").Value = "'2.66"
$ws.Range("E48").Value = "'  +2.43%  "
$ws.Range("D49").Value = "'2.430.60"
$ws.Range("E49").Value = "'  -0.02%  "
$ws.Range("D50").Value = "'69.56"
$ws.Range("E50").Value = "'  -1.68%  "
$ws.Range("B51").Value = "'Stacks"
$ws.Range("C51").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.45"
$ws.Range("E51").Value = "'  +3.92%  "
